# New PO forecast model
# Refreshes the weekly quantity, monthly trend, and PO forecast sheets with
# a new week of actuals/history and an updated forward-looking forecast.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Weekly Quantity" - append the newly observed week
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$lastRow1 = $ws1.UsedRange.Rows.Count

$ws1.Cells.Item($lastRow1 + 1, 1).Value = 45662.99999999999
$ws1.Cells.Item($lastRow1 + 1, 1).NumberFormat = $ws1.Cells.Item($lastRow1, 1).NumberFormat
$ws1.Cells.Item($lastRow1 + 1, 2).Value = 103

# ---------------------------------------------------------------------
# Sheet 2: "Monthly Trend" - append the newly observed month
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Monthly Trend")
$lastRow2 = $ws2.UsedRange.Rows.Count

$ws2.Cells.Item($lastRow2 + 1, 1).Value = 45688.99999999999
$ws2.Cells.Item($lastRow2 + 1, 1).NumberFormat = $ws2.Cells.Item($lastRow2, 1).NumberFormat
$ws2.Cells.Item($lastRow2 + 1, 2).Value = 103

# ---------------------------------------------------------------------
# Sheet 3: "PO Forecast" - refreshed forecast curve
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("PO Forecast")

# Quantity-only updates to existing weekly forecast rows (2-77)
$ws3.Cells.Item(2, 2).Value = 261
$ws3.Cells.Item(3, 2).Value = 260
$ws3.Cells.Item(4, 2).Value = 259
$ws3.Cells.Item(5, 2).Value = 258
$ws3.Cells.Item(6, 2).Value = 257
$ws3.Cells.Item(7, 2).Value = 257
$ws3.Cells.Item(8, 2).Value = 256
$ws3.Cells.Item(9, 2).Value = 255
$ws3.Cells.Item(10, 2).Value = 254
$ws3.Cells.Item(11, 2).Value = 253
$ws3.Cells.Item(12, 2).Value = 253
$ws3.Cells.Item(13, 2).Value = 252
$ws3.Cells.Item(14, 2).Value = 251
$ws3.Cells.Item(15, 2).Value = 250
$ws3.Cells.Item(16, 2).Value = 249
$ws3.Cells.Item(17, 2).Value = 249
$ws3.Cells.Item(18, 2).Value = 248
$ws3.Cells.Item(19, 2).Value = 247
$ws3.Cells.Item(20, 2).Value = 246
$ws3.Cells.Item(21, 2).Value = 245
$ws3.Cells.Item(22, 2).Value = 245
$ws3.Cells.Item(23, 2).Value = 244
$ws3.Cells.Item(24, 2).Value = 243
$ws3.Cells.Item(25, 2).Value = 242
$ws3.Cells.Item(27, 2).Value = 241
$ws3.Cells.Item(28, 2).Value = 240
$ws3.Cells.Item(29, 2).Value = 239
$ws3.Cells.Item(32, 2).Value = 237
$ws3.Cells.Item(41, 2).Value = 229
$ws3.Cells.Item(43, 2).Value = 225
$ws3.Cells.Item(45, 2).Value = 224
$ws3.Cells.Item(46, 2).Value = 223
$ws3.Cells.Item(47, 2).Value = 222
$ws3.Cells.Item(48, 2).Value = 221
$ws3.Cells.Item(49, 2).Value = 219
$ws3.Cells.Item(50, 2).Value = 218
$ws3.Cells.Item(51, 2).Value = 217
$ws3.Cells.Item(52, 2).Value = 217
$ws3.Cells.Item(53, 2).Value = 216
$ws3.Cells.Item(54, 2).Value = 215
$ws3.Cells.Item(55, 2).Value = 214
$ws3.Cells.Item(56, 2).Value = 213
$ws3.Cells.Item(57, 2).Value = 213
$ws3.Cells.Item(58, 2).Value = 212
$ws3.Cells.Item(59, 2).Value = 211
$ws3.Cells.Item(60, 2).Value = 209
$ws3.Cells.Item(61, 2).Value = 209
$ws3.Cells.Item(62, 2).Value = 208
$ws3.Cells.Item(63, 2).Value = 207
$ws3.Cells.Item(64, 2).Value = 206
$ws3.Cells.Item(65, 2).Value = 205
$ws3.Cells.Item(66, 2).Value = 205
$ws3.Cells.Item(67, 2).Value = 204
$ws3.Cells.Item(68, 2).Value = 203
$ws3.Cells.Item(69, 2).Value = 198
$ws3.Cells.Item(70, 2).Value = 197
$ws3.Cells.Item(71, 2).Value = 197
$ws3.Cells.Item(72, 2).Value = 194
$ws3.Cells.Item(73, 2).Value = 193
$ws3.Cells.Item(74, 2).Value = 193
$ws3.Cells.Item(75, 2).Value = 192
$ws3.Cells.Item(76, 2).Value = 183
$ws3.Cells.Item(77, 2).Value = 181

# Forecast rows 78-85 roll forward to new future weeks with new quantities
$ws3.Cells.Item(78, 1).Value = 45662.99999999999
$ws3.Cells.Item(78, 2).Value = 177
$ws3.Cells.Item(79, 1).Value = 45669.99999999999
$ws3.Cells.Item(79, 2).Value = 177
$ws3.Cells.Item(80, 1).Value = 45676.99999999999
$ws3.Cells.Item(80, 2).Value = 176
$ws3.Cells.Item(81, 1).Value = 45683.99999999999
$ws3.Cells.Item(81, 2).Value = 175
$ws3.Cells.Item(82, 1).Value = 45690.99999999999
$ws3.Cells.Item(82, 2).Value = 174
$ws3.Cells.Item(83, 1).Value = 45697.99999999999
$ws3.Cells.Item(83, 2).Value = 173
$ws3.Cells.Item(84, 1).Value = 45704.99999999999
$ws3.Cells.Item(84, 2).Value = 173
$ws3.Cells.Item(85, 1).Value = 45711.99999999999
$ws3.Cells.Item(85, 2).Value = 172

# New forecast week appended at the end of the series
$ws3.Cells.Item(86, 1).Value = 45718.99999999999
$ws3.Cells.Item(86, 1).NumberFormat = $ws3.Cells.Item(85, 1).NumberFormat
$ws3.Cells.Item(86, 2).Value = 171
